$d = $word.ActiveDocument

# The document starts as a single empty paragraph that is centered and
# whose (paragraph-mark) formatting is bold Times New Roman 12pt. The
# edit removes the centering, fills the paragraph with a run of plain
# (non-bold) Times New Roman 12pt body text, and un-bolds the paragraph
# mark formatting to match.

$para = $d.Paragraphs.First

# Remove the centered alignment (back to the document default).
$para.Format.Alignment = 0

$text = "The progressive development of man is vitally dependent on invention. It is the most important product of his creative brain. Its ultimate purpose is the complete mastery of mind over the material world, the harnessing of the forces of nature to human needs. This is the difficult task of the inventor who is often misunderstood and unrewarded. But he finds ample compensation in the pleasing exercises of his powers and in the knowledge of being one of that exceptionally privileged class without whom the race would have long ago perished in the bitter struggle against pitiless elements."

$para.Range.InsertBefore($text)

# Re-grab the (now text-bearing) paragraph range and normalize the font:
# Times New Roman, 12pt, not bold -- for both the new run and the
# paragraph mark that trails it.
$rng = $d.Paragraphs.First.Range
$rng.Font.NameAscii = "Times New Roman"
$rng.Font.NameFarEast = "Times New Roman"
$rng.Font.NameOther = "Times New Roman"
$rng.Font.NameBi = "Times New Roman"
$rng.Font.Size = 12
$rng.Font.SizeBi = 12
$rng.Font.Bold = 0

Write-Host "done"
